$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 01:52"
$ws.Range("B4").Value = 792076
$ws.Range("C4").Value = 27440
$ws.Range("D4").Value = 71947
$ws.Range("E4").Value = 677646
$ws.Range("F4").Value = 13946
$ws.Range("G4").Value = 1908
$ws.Range("H4").Value = 42483
$ws.Range("B16").Value = 36829
$ws.Range("C16").Value = 1773
$ws.Range("D16").Value = 12586
$ws.Range("E16").Value = 22553
$ws.Range("G16").Value = 103
$ws.Range("H16").Value = 1690
$ws.Range("A56").Value = "Argentina"
$ws.Range("B56").Value = 3031
$ws.Range("C56").Value = 90
$ws.Range("D56").Value = 737
$ws.Range("E56").Value = 2152
$ws.Range("F56").Value = 123
$ws.Range("G56").Value = 8
$ws.Range("H56").Value = 142
$ws.Range("A57").Value = "Banglades"
$ws.Range("B57").Value = 2948
$ws.Range("C57").Value = 492
$ws.Range("D57").Value = 85
$ws.Range("E57").Value = 2762
$ws.Range("F57").Value = 1
$ws.Range("G57").Value = 10
$ws.Range("H57").Value = 101
$ws.Range("B87").Value = 879
$ws.Range("C87").Value = 32
$ws.Range("D87").Value = 287
$ws.Range("E87").Value = 582
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 10
$ws.Range("A97").Value = "Guinea"
$ws.Range("B97").Value = 622
$ws.Range("C97").Value = 43
$ws.Range("D97").Value = 122
$ws.Range("E97").Value = 495
$ws.Range("F97").Value = 0
$ws.Range("H97").Value = 5
$ws.Range("A98").Value = "Albania"
$ws.Range("B98").Value = 584
$ws.Range("C98").Value = 22
$ws.Range("D98").Value = 327
$ws.Range("E98").Value = 231
$ws.Range("F98").Value = 5
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 26
$ws.Range("A99").Value = "Burkina Faso"
$ws.Range("B99").Value = 581
$ws.Range("C99").Value = 5
$ws.Range("D99").Value = 357
$ws.Range("E99").Value = 186
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 38
$ws.Range("B102").Value = 535
$ws.Range("C102").Value = 18
$ws.Range("D102").Value = 313
$ws.Range("E102").Value = 212
$ws.Range("B138").Value = 119
$ws.Range("C138").Value = 8
$ws.Range("E138").Value = 107
$ws.Range("D139").Value = 22
$ws.Range("E139").Value = 84
$ws.Range("D146").Value = 37
$ws.Range("E146").Value = 44
$ws.Range("D150").Value = 23
$ws.Range("E150").Value = 47
$ws.Range("A158").Value = "Haiti"
$ws.Range("B158").Value = 57
$ws.Range("C158").Value = 10
$ws.Range("D158").Value = 0
$ws.Range("E158").Value = 54
$ws.Range("F158").Value = 0
$ws.Range("H158").Value = 3
$ws.Range("A159").Value = "Polinesia Francesa"
$ws.Range("D159").Value = 19
$ws.Range("E159").Value = 37
$ws.Range("F159").Value = 1
$ws.Range("A160").Value = "Uganda"
$ws.Range("B160").Value = 56
$ws.Range("C160").Value = 1
$ws.Range("D160").Value = 38
$ws.Range("E160").Value = 18
$ws.Range("H160").Value = 0
$ws.Range("A161").Value = "Benin"
$ws.Range("B161").Value = 54
$ws.Range("C161").Value = 19
$ws.Range("D161").Value = 27
$ws.Range("E161").Value = 26
$ws.Range("A162").Value = "Libia"
$ws.Range("B162").Value = 51
$ws.Range("D162").Value = 15
$ws.Range("E162").Value = 35
$ws.Range("H162").Value = 1
$ws.Range("A163").Value = "Guinea-Bisau"
$ws.Range("B163").Value = 50
$ws.Range("D163").Value = 3
$ws.Range("E163").Value = 47
$ws.Range("H163").Value = 0
$ws.Range("A166").Value = "Eritrea"
$ws.Range("D166").Value = 3
$ws.Range("H166").Value = 0
$ws.Range("A167").Value = "Puerto Rico"
$ws.Range("D167").Value = 1
$ws.Range("H167").Value = 2
$ws.Range("A212").Value = "Sudan del Sur"
$ws.Range("A213").Value = "Santo Tome y Principe"
$ws.Range("A215").Value = "Yemen"
$ws.Range("A216").Value = "San Pedro y Miquelon"
